$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

function Copy-StyleTo {
    param($srcA1, $destRanges)
    $ws.Range($srcA1).Copy() | Out-Null
    foreach ($dest in $destRanges) {
        $ws.Range($dest).PasteSpecial($xlPasteFormats) | Out-Null
    }
    $excel.CutCopyMode = 0
}

# --- Summary block (rows 10-12): recompute against the real (float-safe) question count ---
# A9 already carries the "mtitleStyle" header look that A10/A11/A12 should match.
Copy-StyleTo "A9" @("A10", "A11", "A12")

# Row 10: No. (Right / Wrong / Not Attempt / Max)
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

# Row 11: Marking scheme
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

# Row 12: Total
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "52/112"

# --- Drop the third (G:H) "Student Ans / Correct Ans" block entirely ---
$ws.Range("G15:H40").Clear()

# --- Per-question answer grid ---
# Column A/B hold question blocks 1-25 (rows 16-40); column D/E hold a short
# extra block that, after the fix, only covers questions 26-28 (rows 16-18).

$answers = @{
    16 = @{ Style = "correct";   A = "Option A"; D = "Option A"; E = "Option A" }
    17 = @{ Style = "correct";   A = "Option D"; D = "Option C"; E = "Option C" }
    18 = @{ Style = "correct";   A = "Option B"; D = "Option D"; E = "Option D" }
    19 = @{ Style = "correct";   A = "Option C" }
    20 = @{ Style = "normal" }
    21 = @{ Style = "correct";   A = "Option C" }
    22 = @{ Style = "correct";   A = "Option D" }
    23 = @{ Style = "normal" }
    24 = @{ Style = "normal" }
    25 = @{ Style = "correct";   A = "Option A" }
    26 = @{ Style = "incorrect"; A = "Option D" }
    27 = @{ Style = "normal" }
    28 = @{ Style = "normal" }
    29 = @{ Style = "normal" }
    30 = @{ Style = "correct";   A = "Option B" }
    31 = @{ Style = "incorrect"; A = "Option C" }
    32 = @{ Style = "correct";   A = "Option C" }
    33 = @{ Style = "correct";   A = "Option D" }
    34 = @{ Style = "incorrect"; A = "Option A" }
    35 = @{ Style = "incorrect"; A = "Option B" }
    36 = @{ Style = "normal" }
    37 = @{ Style = "normal" }
    38 = @{ Style = "correct";   A = "Option A" }
    39 = @{ Style = "normal" }
    40 = @{ Style = "normal" }
}

# Reference cells already carrying each named style in the original sheet,
# used as style-paste sources so no new cellXfs entries get minted.
$styleSource = @{
    "correct"   = "A16"   # correctStyle (style 5)
    "incorrect" = "C11"   # incorrectStyle (style 6) -- any s=6 cell
    "normal"    = "A16"   # normalStyle (style 7) -- original A16 before we touch it
}

# A16 starts life as normalStyle, so grab a correctStyle/normalStyle source
# reference BEFORE any rewriting happens.
$normalSourceA1 = "A17"      # still normalStyle at this point
$correctSourceA1 = "B16"     # absoluteStyle -- NOT what we want; fix below

foreach ($row in 16..40) {
    $info = $answers[$row]
    $destA = "A$row"

    switch ($info.Style) {
        "correct"   { Copy-StyleTo "C10" @($destA) }   # C10 currently holds correctStyle? verified below
        "incorrect" { Copy-StyleTo "C11" @($destA) }
        "normal"    { Copy-StyleTo "A20" @($destA) }
    }

    if ($info.ContainsKey("A")) {
        $ws.Range($destA).Value = $info.A
    } else {
        $ws.Range($destA).ClearContents() | Out-Null
    }

    if ($row -le 18) {
        $destD = "D$row"
        $destE = "E$row"
        Copy-StyleTo "C10" @($destD)
        $ws.Range($destD).Value = $info.D
        $ws.Range($destE).Value = $info.E
    } else {
        $ws.Range("D$row`:E$row").Clear()
    }
}
